{"js": "// The document contains five `<id>...</id>` tags, each originally split across\n// three runs: \"<id>\", \"p098v_aN\", \"</id>\". This edit collapses each trio into a\n// single run (taking on the Courier-New \"tag\" formatting) and renumbers the\n// identifier from \"p098v_aN\" to \"p098v_N\" (dropping the \"a\").\nconst body = context.document.body;\n\nfor (let n = 1; n <= 5; n++) {\n  const oldTag = `<id>p098v_a${n}</id>`;\n  const newTag = `<id>p098v_${n}</id>`;\n\n  const results = body.search(oldTag, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  // Replacing the whole found range (which spans the original 3 runs) with the\n  // new text merges it into a single run using the first run's formatting.\n  results.items[0].insertText(newTag, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The document contains five `<id>...</id>` tags, each originally split across\n# three runs: \"<id>\", \"p098v_aN\", \"</id>\". This edit collapses each trio into a\n# single run (taking on the Courier-New \"tag\" formatting) and renumbers the\n# identifier from \"p098v_aN\" to \"p098v_N\" (dropping the \"a\").\n$d = $word.ActiveDocument\n\nfor ($n = 1; $n -le 5; $n++) {\n    $old = \"<id>p098v_a$n</id>\"\n    $new = \"<id>p098v_$n</id>\"\n\n    $range = $d.Content\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #              Format, ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll) \u2014 only one match\n    # exists per tag, so \"replace all\" just replaces that single hit. The\n    # whole matched range (spanning the original 3 runs) is replaced with the\n    # new text, which merges it into a single run using the first run's\n    # formatting.\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
